$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $d.Content.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2) | Out-Null
}

# Title (appears twice: main heading and bold text near the end)
Replace-Text "Play Great Book of Magic Deluxe free: Review and Pros & Cons" "Play Great Book of Magic Deluxe for Free"

# "What we like" bullet list
Replace-Text "Exciting special features including Scatter/Wild symbol and Bonus symbol during Free Spins" "Exciting special features and tools"
Replace-Text "Unique Wazdan Features allow players to personalize the gaming experience" "Customizable gameplay options"
Replace-Text "Customizable gameplay with Volatility Levels, Ultra Lite mode, Ultra Fast mode, and Big Screen mode" "Immersive visuals and atmosphere"
Replace-Text "Great graphics and soundtrack create a magical atmosphere" "Potential for big wins"

# "What we don't like" bullet list
Replace-Text "RTP is good but average for an online slot game" "Limited special features"
Replace-Text "Game may not appeal to those who are not fans of the magic/witchcraft theme" "Average RTP"

# Meta description (italic text)
Replace-Text "An in-depth review of Great Book of Magic Deluxe online slot game, including pros and cons. Play for free and experience the magical atmosphere." "Read our review of Great Book of Magic Deluxe to learn more about this exciting online slot game. Play for free and win big!"
